$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "286.22"
Set-TextValue $ws.Range("E2") "1.35%"
Set-TextValue $ws.Range("D3") "29.32"
Set-TextValue $ws.Range("E3") "2.91%"
Set-TextValue $ws.Range("D4") "5.079"
Set-TextValue $ws.Range("E4") "0.48%"
Set-TextValue $ws.Range("D5") "0.06723"
Set-TextValue $ws.Range("E5") "3.25%"
Set-TextValue $ws.Range("D6") "7.320"
Set-TextValue $ws.Range("E6") "1.04%"
Set-TextValue $ws.Range("D7") "1.393"
Set-TextValue $ws.Range("E7") "-1.09%"
Set-TextValue $ws.Range("D8") "0.9005"
Set-TextValue $ws.Range("E8") "-2.36%"
Set-TextValue $ws.Range("E9") "2.29%"
Set-TextValue $ws.Range("E10") "7.24%"
Set-TextValue $ws.Range("D11") "0.07569"
Set-TextValue $ws.Range("E11") "-0.69%"
Set-TextValue $ws.Range("D12") "0.02924"
Set-TextValue $ws.Range("E12") "0.77%"
Set-TextValue $ws.Range("D13") "0.08996"
Set-TextValue $ws.Range("E13") "0.60%"
Set-TextValue $ws.Range("D14") "0.001593"
Set-TextValue $ws.Range("E14") "0.34%"
Set-TextValue $ws.Range("D15") "0.04484"
Set-TextValue $ws.Range("E15") "1.43%"
Set-TextValue $ws.Range("D16") "0.0006452"
Set-TextValue $ws.Range("E16") "0.86%"
Set-TextValue $ws.Range("D17") "0.006493"
Set-TextValue $ws.Range("E17") "7.31%"
Set-TextValue $ws.Range("D18") "3.449"
Set-TextValue $ws.Range("E18") "0.32%"
Set-TextValue $ws.Range("D19") "3.437"
Set-TextValue $ws.Range("E19") "1.49%"
Set-TextValue $ws.Range("E20") "-0.49%"
Set-TextValue $ws.Range("E21") "0.51%"
Set-TextValue $ws.Range("E22") "2.34%"
Set-TextValue $ws.Range("D23") "4.048"
Set-TextValue $ws.Range("E23") "0.55%"
Set-TextValue $ws.Range("D24") "0.1581"
Set-TextValue $ws.Range("E24") "2.12%"
Set-TextValue $ws.Range("D25") "0.001200"
Set-TextValue $ws.Range("E25") "1.03%"
Set-TextValue $ws.Range("D26") "0.004370"
Set-TextValue $ws.Range("E26") "-0.63%"
Set-TextValue $ws.Range("E27") "-6.89%"
Set-TextValue $ws.Range("D28") "0.0001617"
Set-TextValue $ws.Range("E28") "-0.25%"
Set-TextValue $ws.Range("D40") "0.04240"
Set-TextValue $ws.Range("E40") "2.20%"
Set-TextValue $ws.Range("D41") "0.006769"
Set-TextValue $ws.Range("E41") "1.18%"
Set-TextValue $ws.Range("E42") "1.34%"
Set-TextValue $ws.Range("D43") "0.002186"
Set-TextValue $ws.Range("E43") "1.80%"
Set-TextValue $ws.Range("D44") "0.01148"
Set-TextValue $ws.Range("E44") "-5.65%"
Set-TextValue $ws.Range("D45") "0.00005721"
Set-TextValue $ws.Range("E45") "1.58%"
Set-TextValue $ws.Range("D46") "1.937"
Set-TextValue $ws.Range("E46") "-1.45%"
Set-TextValue $ws.Range("E47") "15.05%"
